$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Fill in the previously empty "Name" value cell with the generated type name
$ws.Range("B4").Value = "StatutetatcivilVs"

# Update the Date value cell to reflect the new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
